$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's data point needs to land in row 2 (most recent first), so push
# every existing data row (2..21) down by one (to 3..22), working from the
# bottom up so nothing gets clobbered before it's copied.
for ($r = 21; $r -ge 2; $r--) {
    $src = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 20))
    $dst = $ws.Range($ws.Cells.Item($r + 1, 1), $ws.Cells.Item($r + 1, 20))
    $src.Copy($dst)
}

# Populate the now-empty row 2 with this week's record.
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(2, 3).Value = "Los Lagos"
$ws.Cells.Item(2, 4).Value = Get-Date -Year 2022 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100101
$ws.Cells.Item(2, 8).Value = "Berries"
$ws.Cells.Item(2, 9).Value = 100101001
$ws.Cells.Item(2, 10).Value = "Arándano (blue)"
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 7500
$ws.Cells.Item(2, 15).Value = 8000
$ws.Cells.Item(2, 16).Value = 7750
$ws.Cells.Item(2, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(2, 19).Value = 5167
$ws.Cells.Item(2, 20).Value = 1.5
